$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Deialdi ordua: 10:00-11:00" -> "Deialdi ordua: 10:00-13:00" split into
#    three bold runs for the time value.
# ---------------------------------------------------------------------------
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Deialdi ordua:*") {
        $target1 = $p
        break
    }
}

$xml1 = '<w:p ' + $wns + ' w14:paraId="6F4EFF29" w14:textId="7EF35949" w:rsidR="0050238D" w:rsidRDefault="00F51F43" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:keepNext/><w:spacing w:after="80" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Deialdi ordua: </w:t></w:r>' + `
    '<w:r w:rsidR="00120539"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>10:00-1</w:t></w:r>' + `
    '<w:r w:rsidR="00120539"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3</w:t></w:r>' + `
    '<w:r w:rsidR="00120539"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:00</w:t></w:r>' + `
    '</w:p>'
$target1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) " Inplementazioa ez dakit zen" -> " Inplementazioarekin ez dakit zenbat
#    denbora beharko dudan." split across four runs.
# ---------------------------------------------------------------------------
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Inplementazioa ez dakit zen*") {
        $target2 = $p
        break
    }
}

$xml2 = '<w:p ' + $wns + ' w14:paraId="0E807601" w14:textId="55A3E24F" w:rsidR="00382D35" w:rsidRDefault="00382D35" w:rsidP="00382D35" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:keepNext/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="28"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Azkar ari naiz aurreratzen baina oraindik falta da. </w:t></w:r>' + `
    '<w:r w:rsidR="006B7EC8"><w:t>Memoria</w:t></w:r>' + `
    '<w:r w:rsidR="00016954"><w:t xml:space="preserve">, eranskinak, aurkezpena eta posterra </w:t></w:r>' + `
    '<w:r w:rsidR="006B7EC8"><w:t>bukatzeko denbora izango dut nire ustez.</w:t></w:r>' + `
    '<w:r w:rsidR="00016954"><w:t xml:space="preserve"> Inplementazioa</w:t></w:r>' + `
    '<w:r w:rsidR="00016954"><w:t>rekin</w:t></w:r>' + `
    '<w:r w:rsidR="00016954"><w:t xml:space="preserve"> ez dakit zen</w:t></w:r>' + `
    '<w:r w:rsidR="00016954"><w:t>bat denbora beharko dudan.</w:t></w:r>' + `
    '</w:p>'
$target2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) "Segurunez" (with proofErr spell-check wrapper) -> "Seguru" + "e" + "nez"
#    (no proofErr), leading space merged into the first new run.
# ---------------------------------------------------------------------------
$target3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Segurunez*") {
        $target3 = $p
        break
    }
}

$xml3 = '<w:p ' + $wns + ' w14:paraId="4E7CF090" w14:textId="71ECE315" w:rsidR="00382D35" w:rsidRDefault="00382D35" w:rsidP="00382D35" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:keepNext/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="28"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr>' + `
    '<w:r><w:t>Hurrengo asterako</w:t></w:r>' + `
    '<w:r w:rsidR="00673949"><w:t xml:space="preserve"> erabaki behar da defentsa uztailean edo irailean egin. </w:t></w:r>' + `
    '<w:r w:rsidR="006B7EC8"><w:t>Bilera eguna azkena denez orduan erabaki dezakegu.</w:t></w:r>' + `
    '<w:r w:rsidR="008E240F"><w:t xml:space="preserve"> Seguru</w:t></w:r>' + `
    '<w:r w:rsidR="008E240F"><w:t>e</w:t></w:r>' + `
    '<w:r w:rsidR="008E240F"><w:t xml:space="preserve">nez ez du denborarik emango uztailean </w:t></w:r>' + `
    '<w:r w:rsidR="008E240F"><w:lastRenderedPageBreak/><w:t>aurkezteko.</w:t></w:r>' + `
    '</w:p>'
$target3.Range.InsertXML($xml3)

Write-Output "done"
